$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.30448347740663451
$ws.Cells.Item(2, 1).Value = -0.0059999999712410101
$ws.Cells.Item(3, 1).Value = -0.0039999999712607703
$ws.Cells.Item(4, 1).Value = -0.0079999999509698938
$ws.Cells.Item(5, 1).Value = -0.0029999999712391201
$ws.Cells.Item(6, 1).Value = -0.0019999999710869076
$ws.Cells.Item(7, 1).Value = -0.0099999999351458868
$ws.Cells.Item(8, 1).Value = -0.0099999999349331681
$ws.Cells.Item(9, 1).Value = -0.0019999999715545336
$ws.Cells.Item(10, 1).Value = -0.0019999999723090411
$ws.Cells.Item(11, 1).Value = 0.032489237710034757
$ws.Cells.Item(12, 1).Value = -0.0034999999649802938
$ws.Cells.Item(13, 1).Value = -0.0034999999628313461
$ws.Cells.Item(14, 1).Value = -0.0079999999420135026
$ws.Cells.Item(15, 1).Value = -0.00099999997319777378
$ws.Cells.Item(16, 1).Value = -0.0019999999683135705
$ws.Cells.Item(17, 1).Value = -0.0019999999678503855
$ws.Cells.Item(18, 1).Value = -0.0039999999586770585
$ws.Cells.Item(19, 1).Value = -0.050977678753525879
$ws.Cells.Item(20, 1).Value = -0.0039999999781414886
$ws.Cells.Item(21, 1).Value = -0.0039999999779229967
$ws.Cells.Item(22, 1).Value = -0.0039999999777498019
$ws.Cells.Item(23, 1).Value = -0.034435063589466353
$ws.Cells.Item(24, 1).Value = -0.019999999887263975
$ws.Cells.Item(25, 1).Value = -0.079235648792076674
$ws.Cells.Item(26, 1).Value = -0.0024999999680730411
$ws.Cells.Item(27, 1).Value = -0.0024999999668398054
$ws.Cells.Item(28, 1).Value = -0.0019999999642825728
$ws.Cells.Item(29, 1).Value = -0.0069999999384222633
$ws.Cells.Item(30, 1).Value = -0.059999999698133966
$ws.Cells.Item(31, 1).Value = 0.028507653454292736
$ws.Cells.Item(32, 1).Value = -0.0099999999240729665
$ws.Cells.Item(33, 1).Value = -0.003999999950998756
